# this commit for new github
#
# Applies the tracked changes:
#   1. Refresh the "automatic date" placeholder text (2/21/2021 -> 3/5/2021)
#      on the slide master and every slide layout.
#   2. Add body text to the (previously empty) subtitle placeholder on
#      slide 1: two runs, "I changed little bit " + "in this file".

function Set-DatePlaceholderText {
    param($shapes, $text)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $isDate = $true
            }
        } catch {
        }
        if ($isDate) {
            $sh.TextFrame.TextRange.Text = $text
        }
    }
}

$p = $ppt.ActivePresentation

# --- 1. Update the cached "datetimeFigureOut" placeholder text everywhere ---
$newDate = "3/5/2021"

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes $newDate

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# --- 2. Fill in the subtitle placeholder text on slide 1 ---
$slide = $p.Slides.Item(1)
$subtitle = $slide.Shapes.Item(2)
$tr = $subtitle.TextFrame.TextRange
$tr.Text = "I changed little bit "
$tr.InsertAfter("in this file")
